# QA Compiler: Fix Actual Issues % negative value bug
#
# The underlying bug (clamping actual_pct to 0-100%) lives in the Python
# report generator, not in this workbook. The committed fix simply
# re-renders the already-100%-capped "Actual Issues" percentage cells
# from "100.0%" to "100%" on the DAILY and TOTAL sheets (the hidden
# _DAILY_DATA sheet and the "Comp %" columns are untouched).
#
# Every affected cell stores its percentage as literal text (t="inlineStr"
# in the source file), not as a numeric percent. Assigning a percent-
# looking string straight to Range.Value (e.g. "100%") gets auto-parsed
# into a numeric 1.0 with a new percentage number format / style, which
# would NOT match the original formatting. To keep the literal text and
# the cell's existing style untouched, we stage the text in a scratch
# formula cell (="100%") so it is produced as a string, copy it, and
# paste-special *values only* into each target cell - this swaps just the
# displayed text without disturbing the target's style index.

$wb = $excel.ActiveWorkbook

function Set-ClampedPercentText {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Text
    )

    # Scratch cell far outside the used range of either sheet.
    $scratch = $Worksheet.Range("ZZ1")
    $scratch.Formula = '="' + $Text + '"'

    $target = $Worksheet.Range($Address)
    $scratch.Copy()
    $target.PasteSpecial(-4163)   # xlPasteValues - value only, keep target's own style

    $scratch.Clear()
}

# --- DAILY sheet ("Actual Issues" columns: E = Alice, I = Bob, M = John) ---
$daily = $wb.Worksheets.Item("DAILY")
foreach ($addr in @("E5", "I6", "M6", "E8", "I8", "M8")) {
    Set-ClampedPercentText $daily $addr "100%"
}

# --- TOTAL sheet ("Actual Issues" column: C) ---
$total = $wb.Worksheets.Item("TOTAL")
foreach ($addr in @("C3", "C4", "C5", "C6", "C9")) {
    Set-ClampedPercentText $total $addr "100%"
}
